$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on Price (D) and Volume(1h) (E) columns so that
# numeric-looking strings (e.g. "1.000", "5.110") are preserved exactly
# as text and not coerced into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.761.50'
$ws.Range("E2").Value = '  +1.19%  '
$ws.Range("D3").Value = '1.885.18'
$ws.Range("E3").Value = '  +1.99%  '
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '238.84'
$ws.Range("E5").Value = '  +2.30%  '
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = '0.4765'
$ws.Range("E7").Value = '  +2.02%  '
$ws.Range("D8").Value = '0.2862'
$ws.Range("E8").Value = '  +4.74%  '
$ws.Range("D9").Value = '0.06563'
$ws.Range("E9").Value = '  +4.28%  '
$ws.Range("D10").Value = '18.90'
$ws.Range("E10").Value = '  +16.26%  '
$ws.Range("D11").Value = '1.878.06'
$ws.Range("E11").Value = '  +1.70%  '
$ws.Range("B12").Value = 'Litecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D12").Value = '96.26'
$ws.Range("E12").Value = '  +14.72%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.07563'
$ws.Range("E13").Value = '  +1.41%  '
$ws.Range("D14").Value = '5.110'
$ws.Range("E14").Value = '  +3.44%  '
$ws.Range("D15").Value = '0.6549'
$ws.Range("E15").Value = '  +5.61%  '
$ws.Range("D16").Value = '304.42'
$ws.Range("E16").Value = '  +32.54%  '
$ws.Range("D17").Value = '30.759.24'
$ws.Range("E17").Value = '  +1.38%  '
$ws.Range("D18").Value = '13.17'
$ws.Range("E18").Value = '  +6.46%  '
$ws.Range("D19").Value = '0.9992'
$ws.Range("D20").Value = '0.000007579'
$ws.Range("E20").Value = '  +3.31%  '
$ws.Range("D21").Value = '2.124.30'
$ws.Range("E21").Value = '  +2.58%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = '5.112'
$ws.Range("E23").Value = '  +3.85%  '
$ws.Range("D24").Value = '6.172'
$ws.Range("E24").Value = '  +5.25%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '9.278'
$ws.Range("E25").Value = '  +1.50%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '166.99'
$ws.Range("E26").Value = '  +0.90%  '
$ws.Range("D27").Value = '20.26'
$ws.Range("E27").Value = '  +13.57%  '
$ws.Range("D28").Value = '1.951'
$ws.Range("E28").Value = '  +4.30%  '
$ws.Range("D29").Value = '0.1080'
$ws.Range("E29").Value = '  +5.52%  '
$ws.Range("D30").Value = '1.352'
$ws.Range("E30").Value = '  -1.51%  '
$ws.Range("D31").Value = '4.163'
$ws.Range("E31").Value = '  +1.72%  '
$ws.Range("D32").Value = '3.963'
$ws.Range("E32").Value = '  +3.95%  '
$ws.Range("D33").Value = '0.05052'
$ws.Range("E33").Value = '  +3.70%  '
$ws.Range("D34").Value = '1.173'
$ws.Range("E34").Value = '  +2.68%  '
$ws.Range("D35").Value = '0.7312'
$ws.Range("E35").Value = '  +3.73%  '
$ws.Range("D36").Value = '2.715'
$ws.Range("E36").Value = '  +1.01%  '
$ws.Range("D37").Value = '0.01934'
$ws.Range("E37").Value = '  +1.19%  '
$ws.Range("D38").Value = '2.700'
$ws.Range("E38").Value = '  +1.30%  '
$ws.Range("D39").Value = '2.078'
$ws.Range("E39").Value = '  +7.83%  '
$ws.Range("D40").Value = '0.9030'
$ws.Range("E40").Value = '  +4.38%  '
$ws.Range("D41").Value = '107.62'
$ws.Range("E41").Value = '  +1.78%  '
$ws.Range("D42").Value = '0.9990'
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").Value = '0.4217'
$ws.Range("E43").Value = '  +4.83%  '
$ws.Range("D44").Value = '5.632'
$ws.Range("E44").Value = '  +1.92%  '
$ws.Range("D45").Value = '65.94'
$ws.Range("E45").Value = '  +7.05%  '
$ws.Range("D46").Value = '7.363'
$ws.Range("E46").Value = '  +3.91%  '
$ws.Range("D47").Value = '0.1225'
$ws.Range("E47").Value = '  +1.96%  '
$ws.Range("D48").Value = '8.990'
$ws.Range("E48").Value = '  +4.66%  '
$ws.Range("D49").Value = '34.76'
$ws.Range("E49").Value = '  +4.20%  '
$ws.Range("D50").Value = '0.05610'
$ws.Range("E50").Value = '  +1.70%  '
$ws.Range("D51").Value = '1.390'
$ws.Range("E51").Value = '  +3.24%  '
